# Updated cryptos list on Sun Aug 11 06:59:09 UTC 2024 with GitHub Actions
# Refresh price/volume(1h) figures for each coin row; rows 45/46 (RenderToken,
# EnergySwap) additionally swap places with their refreshed figures.
# Numeric-looking text values are entered with a leading apostrophe so Excel
# keeps them as text (matching the original inlineStr cells) instead of
# auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.050.46'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '2.645.29'
$ws.Range('E3').Value = '  +1.43%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''532.14'
$ws.Range('E5').Value = '  +4.16%  '
$ws.Range('D6').Value = '''155.61'
$ws.Range('E6').Value = '  +0.56%  '
$ws.Range('D8').Value = '''0.592'
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('D9').Value = '''6.63'
$ws.Range('E9').Value = '  -1.00%  '
$ws.Range('E10').Value = '  +4.72%  '
$ws.Range('D11').Value = '''0.351'
$ws.Range('E11').Value = '  +1.52%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '3.110.24'
$ws.Range('E13').Value = '  +1.51%  '
$ws.Range('D14').Value = '61.034.53'
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').Value = '''22.07'
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').Value = '2.654.52'
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = '''354.82'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').Value = '''6.23'
$ws.Range('E21').Value = '  +1.38%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').Value = '''61.70'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('E25').Value = '  +1.56%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Value = '0.0₃0859'
$ws.Range('E27').Value = '  +2.19%  '
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('E30').Value = '  +7.30%  '
$ws.Range('E31').Value = '  +4.13%  '
$ws.Range('D32').Value = '''19.54'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').Value = '''150.01'
$ws.Range('E33').Value = '  -0.71%  '
$ws.Range('D34').Value = '''4.15'
$ws.Range('E34').Value = '  +3.94%  '
$ws.Range('E35').Value = '  +1.33%  '
$ws.Range('D36').Value = '''0.917'
$ws.Range('E36').Value = '  +8.46%  '
$ws.Range('D37').Value = '''0.904'
$ws.Range('E37').Value = '  +1.92%  '
$ws.Range('D38').Value = '''308.33'
$ws.Range('E38').Value = '  +4.50%  '
$ws.Range('E39').Value = '  +1.15%  '
$ws.Range('E40').Value = '  +1.77%  '
$ws.Range('D41').Value = '''0.646'
$ws.Range('E41').Value = '  +3.31%  '
$ws.Range('E42').Value = '  +1.50%  '
$ws.Range('D43').Value = '''0.0562'
$ws.Range('E43').Value = '  +1.15%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''19.87'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '''4.98'
$ws.Range('E46').Value = '  +2.05%  '
$ws.Range('E47').Value = '  +2.63%  '
$ws.Range('D48').Value = '''19.31'
$ws.Range('E48').Value = '  +8.11%  '
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('D50').Value = '1.991.85'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('E51').Value = '  +2.42%  '
